$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.925.48'
$ws.Range("E2").Value = '  -0.93%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.859.01'
$ws.Range("E3").Value = '  -1.37%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.92'
$ws.Range("E5").Value = '  -0.94%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.84'
$ws.Range("E6").Value = '  +0.91%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.857.19'
$ws.Range("E7").Value = '  -1.35%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("E10").Value = '  -0.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.32'
$ws.Range("E11").Value = '  -1.10%  '
$ws.Range("E12").Value = '  -0.52%  '
$ws.Range("E13").Value = '  +1.28%  '
$ws.Range("E14").Value = '  +0.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.505.28'
$ws.Range("E15").Value = '  -1.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.860.95'
$ws.Range("E16").Value = '  -1.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.001.00'
$ws.Range("E17").Value = '  -1.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.21'
$ws.Range("E18").Value = '  +7.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.39'
$ws.Range("E19").Value = '  -0.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.110'
$ws.Range("E20").Value = '  -2.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.76'
$ws.Range("E21").Value = '  -3.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '466.59'
$ws.Range("E22").Value = '  -3.59%  '
$ws.Range("E23").Value = '  +1.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000161'
$ws.Range("E24").Value = '  -4.86%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.34'
$ws.Range("E25").Value = '  -1.24%  '
$ws.Range("E26").Value = '  -0.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.15'
$ws.Range("E27").Value = '  +1.18%  '
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.99'
$ws.Range("E29").Value = '  -1.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.95'
$ws.Range("E30").Value = '  +0.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.009.47'
$ws.Range("E31").Value = '  -1.27%  '
$ws.Range("E32").Value = '  -2.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.31'
$ws.Range("E33").Value = '  -2.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.99'
$ws.Range("E34").Value = '  -2.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.835.13'
$ws.Range("E35").Value = '  -0.54%  '
$ws.Range("E36").Value = '  -2.54%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.02'
$ws.Range("E38").Value = '  -2.10%  '
$ws.Range("B39").Value = 'dogwifhat'
$ws.Range("C39").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.30'
$ws.Range("E39").Value = '  +8.84%  '
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.89'
$ws.Range("E40").Value = '  +0.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("E42").Value = '  -2.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '429.18'
$ws.Range("E43").Value = '  -1.41%  '
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '47.25'
$ws.Range("E46").Value = '  -2.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.54'
$ws.Range("E47").Value = '  +1.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '143.49'
$ws.Range("E48").Value = '  +1.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000271'
$ws.Range("E49").Value = '  +7.20%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.88'
$ws.Range("E50").Value = '  -0.33%  '
$ws.Range("B51").Value = 'Arweave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '39.37'
$ws.Range("E51").Value = '  +0.60%  '
